# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" for the d35d6843-37b0-4f0e-acaa-b4a3859a0c27
# file's row (row 5) on both language sheets, reflecting a fresh handoff just run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$zhcn.Range("D5").Value = "2016-01-28 07:44:56"
$dede.Range("D5").Value = "2016-01-28 07:45:08"
